{"js": "// Highlight the words \"If\", \"If\" and \"Otherwise\" in magenta within the\n// \"Rent DVDs\" use case description paragraphs (migrated Use case diagram\n// note -> marking up the conditional keywords).\n\nconst body = context.document.body;\n\n// --- Edit 1 -------------------------------------------------------------\n// \"...scan the customer card. If the system accepts the card...\"\n// The run \". If the \" becomes \". \" + \"If\" (magenta) + \" the \".\nconst firstIfScope = body.search(\". If the \", { matchCase: true });\nfirstIfScope.load(\"text\");\nawait context.sync();\n\nconst firstIfWord = firstIfScope.items[0].search(\"If\", { matchCase: true });\nawait context.sync();\nfirstIfWord.items[0].font.highlightColor = \"magenta\";\nawait context.sync();\n\n// --- Edit 2 -------------------------------------------------------------\n// \"...stored in the system. If the scanned barcode...\"\n// The existing standalone \"If\" run just gains a magenta highlight.\nconst secondIfScope = body.search(\"system. If the scanned barcode\", { matchCase: true });\nsecondIfScope.load(\"text\");\nawait context.sync();\n\nconst secondIfWord = secondIfScope.items[0].search(\"If\", { matchCase: true });\nawait context.sync();\nsecondIfWord.items[0].font.highlightColor = \"magenta\";\nawait context.sync();\n\n// --- Edit 3 -------------------------------------------------------------\n// \"...the system rejects the DVD. Otherwise, the system will add ...\"\n// Splits into: \"...DVD. \" + \"Otherwise\" (magenta) + \", the system will add \".\nconst otherwiseScope = body.search(\n  \"does not match any stored barcode, the system rejects the DVD. Otherwise, the system will add \",\n  { matchCase: true }\n);\notherwiseScope.load(\"text\");\nawait context.sync();\n\nconst otherwiseWord = otherwiseScope.items[0].search(\"Otherwise\", { matchCase: true });\nawait context.sync();\notherwiseWord.items[0].font.highlightColor = \"magenta\";\nawait context.sync();\n", "ps1": "# Highlight the conditional keywords \"If\", \"If\" and \"Otherwise\" in magenta\n# within the \"Rent DVDs\" use case description (migrated Use case diagram\n# note -> marking up the conditional keywords).\n\n$d = $word.ActiveDocument\n\nfunction Highlight-FirstWord {\n    param($scopeText, $targetWord)\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.Text = $scopeText\n    $find.MatchCase = $true\n    $found = $find.Execute()\n\n    if (-not $found) {\n        throw \"scope text not found: $scopeText\"\n    }\n\n    $offset = $scopeText.IndexOf($targetWord)\n    $wordStart = $range.Start + $offset\n    $wordEnd = $wordStart + $targetWord.Length\n\n    $wordRange = $d.Range($wordStart, $wordEnd)\n    # wdColorIndex 5 = wdPink -> serialises as <w:highlight w:val=\"magenta\"/>\n    $wordRange.Font.HighlightColorIndex = 5\n}\n\n# --- Edit 1 ---------------------------------------------------------------\n# \"...scan the customer card. If the system accepts the card...\"\n# The run \". If the \" becomes \". \" + \"If\" (magenta) + \" the \".\nHighlight-FirstWord \". If the \" \"If\"\n\n# --- Edit 2 -----------------------------------------------------------------\n# \"...stored in the system. If the scanned barcode...\"\n# The existing standalone \"If\" run just gains a magenta highlight.\nHighlight-FirstWord \"system. If the scanned barcode\" \"If\"\n\n# --- Edit 3 -----------------------------------------------------------------\n# \"...the system rejects the DVD. Otherwise, the system will add ...\"\n# Splits into: \"...DVD. \" + \"Otherwise\" (magenta) + \", the system will add \".\nHighlight-FirstWord \"does not match any stored barcode, the system rejects the DVD. Otherwise, the system will add \" \"Otherwise\"\n"}
